$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '24.946.43'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.07%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.701.68'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.95%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.30%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '315.89'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.12%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.002'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.31%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3971'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.33%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4027'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.13%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.469'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.93%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '52.75'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.62%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.002'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.32%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08807'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.64%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '25.94'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.05%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.458'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.36%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.00001350'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.34%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.965'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.17%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.713.17'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.08%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '96.26'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.56%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.07191'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.18%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '20.58'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.85%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.352'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.30%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.002'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.29%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '14.45'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.99%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '24.982.32'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.27%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.351'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.52%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.950'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.29%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '23.75'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +5.31%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.196'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +15.58%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '161.50'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.46%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '149.63'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +8.27%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.331'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.37%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.629'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +29.59%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.906.23'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.34%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08561'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.06%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.03147'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.17%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.048'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.41%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '7.208'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.36%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2862'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.26%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '10.90'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.05%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.09543'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +4.52%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8256'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.27%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '13.97'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.83%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.483'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.13%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '17.27'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.96%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.689'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.49%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.7388'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.11%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.256'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.20%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.407'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.41%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.08742'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +8.31%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.002'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.35%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '139.05'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.07%  '
